# Revert "Adding the RES Hourly Production Forecast to the Portfolio"
# - Shift the Data/Lookup date from 24.09.2024 back to 29.08.2024 (i.e. Column A dates move
#   back by 26 days, and the Column D lookup strings follow the same new date), and restore
#   the original Prediction values in Column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDateText = "29.08.2024"
$dayShift = 26

# New Prediction (Column C) values for rows 2..96 (row index -> value)
$predictionValues = @{
    30 = 0.011; 31 = 0.013; 32 = 0.016; 33 = 0.025; 34 = 0.039; 35 = 0.1;
    36 = 0.126; 37 = 0.112; 38 = 0.134; 39 = 0.156; 40 = 0.164; 41 = 0.149;
    42 = 0.284; 43 = 0.368; 44 = 0.478; 45 = 0.577; 46 = 0.625; 47 = 0.652;
    48 = 0.67;  49 = 0.7;   50 = 0.708; 51 = 0.714; 52 = 0.712; 53 = 0.698;
    54 = 0.6879999999999999; 55 = 0.667; 56 = 0.647; 58 = 0.63; 59 = 0.617;
    60 = 0.594; 61 = 0.569; 62 = 0.539; 63 = 0.501; 64 = 0.487; 65 = 0.456;
    66 = 0.419; 67 = 0.375; 68 = 0.342; 69 = 0.314; 70 = 0.29;  71 = 0.246;
    72 = 0.198; 73 = 0.153; 74 = 0.126; 75 = 0.099; 76 = 0.076; 77 = 0.057;
    78 = 0.042; 79 = 0.03;  80 = 0.022; 81 = 0.015
}

for ($r = 2; $r -le 96; $r++) {
    # Column A: shift the serial date/time back by 26 days, preserving time-of-day fraction
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $oldDate - $dayShift

    # Column B holds the interval number, used to rebuild the Column D lookup text
    $interval = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 4).Value = $newDateText + [string]$interval

    # Column C: apply the restored Prediction value where it changed
    if ($predictionValues.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value = $predictionValues[$r]
    }
}
